$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (shifts the old H:K block to I:L), which also
# carries over formatting from the column to the left (G) for the new
# H cells - exactly mirroring how the "Estado/Transaccion/Fecha/Cuenta"
# block moved one column to the right to make room for "fecha".
$ws.Range("H1:H2").Insert(-4161)

# --- Row 1 (headers) ---
$ws.Range("G1").Value = "cuenta debito"
$ws.Range("H1").Value = "fecha"

# --- Row 2 (data) ---
# Plain text values (not number-like, so Excel/engine keeps them as text
# without any extra type coercion).
$ws.Range("A2").Value = "dmoralesr"
$ws.Range("H2").Value = "R_RENEWAL +"
$ws.Range("I2").Value = "PASSED"
$ws.Range("J2").Value = "AAACT231773NFXKK3"
$ws.Range("K2").Value = "26 jun. 2023, 18:00:45"

# Number-like text values that must stay text (quote-prefix) inside the
# B:H data block, matching the sibling cells in that block.
$ws.Range("C2").Value = "'4862917"
$ws.Range("E2").Value = "'5000"
$ws.Range("G2").Value = "'1010506979"

# L2 is number-like text too, but (like the other cells in the I:L
# block) carries no quote-prefix styling. Stage it with a quoted value
# in a scratch cell, then copy only the VALUE (not the formatting) over,
# so it lands as plain shared-string text with the default style.
$ws.Range("Z1").Value = "'1010825578"
$ws.Range("Z1").Copy()
$ws.Range("L2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Recompute best-fit column widths to match the new content.
$ws.Columns.AutoFit()

# Selection as left by the editor.
$ws.Range("H2").Select()
